$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 18: RR vs KKR results entered for row 27
$ws.Range("E27").Value = 40
$ws.Range("H27").Value = 60
$ws.Range("K27").Value = 20
$ws.Range("N27").Value = 100
$ws.Range("Q27").Value = 80
$ws.Range("T27").Value = 0
